$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.434.79'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '2.379.16'
$ws.Range("E3").Value = '  +4.91%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'235.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").Value = "'0.649"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("D7").Value = "'71.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +11.93%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").Value = "'56.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").Value = "'27.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '2.724.50'
$ws.Range("E13").Value = '  +4.70%  '
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = "'16.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = "'6.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.18%  '
$ws.Range("D17").Value = "'0.854"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").Value = '2.371.42'
$ws.Range("E18").Value = '  +4.67%  '
$ws.Range("D19").Value = '43.425.89'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'6.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.08%  '
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").Value = "'74.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").Value = "'250.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = "'3.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.41%  '
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = "'2.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.86%  '
$ws.Range("D30").Value = "'174.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  +5.02%  '
$ws.Range("E32").Value = '  -6.28%  '
$ws.Range("D33").Value = "'0.127"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = "'5.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("E36").Value = '  +1.84%  '
$ws.Range("E37").Value = '  +6.89%  '
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("D39").Value = "'3.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("D40").Value = "'0.0257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").Value = "'8.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = "'18.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.25%  '
$ws.Range("E44").Value = '  +7.88%  '
$ws.Range("D45").Value = "'100.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").Value = "'4.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("E47").Value = '  +2.10%  '
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E49").Value = '  -7.33%  '
$ws.Range("D50").Value = '1.445.79'
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").Value = '2.601.13'
$ws.Range("E51").Value = '  +5.04%  '
